$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.4946373333333334
$ws.Range("H2").Value = 1.483912
$ws.Range("I2").Value = 0.1240056293871995
$ws.Range("J2").Value = 0.1240056293871996
$ws.Range("M2").Value = 0.02123233333333334
$ws.Range("N2").Value = 0.063697
$ws.Range("O2").Value = 0.08772724768620539
$ws.Range("P2").Value = 0.08772724768620538
$ws.Range("Q2").Value = 0.01050230474044445
$ws.Range("R2").Value = 0.09452074266400001
$ws.Range("S2").Value = 0.01087867256373464
$ws.Range("T2").Value = 0.01087867256373464

# Row 3
$ws.Range("G3").Value = 0.4946373333333334
$ws.Range("H3").Value = 1.483912
$ws.Range("I3").Value = 0.1240056293871995
$ws.Range("J3").Value = 0.1240056293871996
$ws.Range("N3").Value = 0.6623830000000001
$ws.Range("O3").Value = 0.9122727523137947
$ws.Range("P3").Value = 0.9122727523137947
$ws.Range("Q3").Value = 0.1092131202551111
$ws.Range("R3").Value = 0.9829180822960002
$ws.Range("S3").Value = 0.1131269568234649
$ws.Range("T3").Value = 0.1131269568234649

# Row 4
$ws.Range("I4").Value = 0.3425317985918844
$ws.Range("J4").Value = 0.3425317985918844
$ws.Range("M4").Value = 0.02123233333333334
$ws.Range("N4").Value = 0.063697
$ws.Range("O4").Value = 0.08772724768620539
$ws.Range("P4").Value = 0.08772724768620538
$ws.Range("Q4").Value = 0.02900975826566667
$ws.Range("R4").Value = 0.261087824391
$ws.Range("S4").Value = 0.03004937193547166
$ws.Range("T4").Value = 0.03004937193547166

# Row 5
$ws.Range("I5").Value = 0.3425317985918844
$ws.Range("J5").Value = 0.3425317985918844
$ws.Range("N5").Value = 0.6623830000000001
$ws.Range("O5").Value = 0.9122727523137947
$ws.Range("P5").Value = 0.9122727523137947
$ws.Range("S5").Value = 0.3124824266564127
$ws.Range("T5").Value = 0.3124824266564128

# Row 6
$ws.Range("G6").Value = 2.127891333333333
$ws.Range("H6").Value = 6.383674
$ws.Range("I6").Value = 0.5334625720209161
$ws.Range("J6").Value = 0.5334625720209161
$ws.Range("M6").Value = 0.02123233333333334
$ws.Range("N6").Value = 0.063697
$ws.Range("O6").Value = 0.08772724768620539
$ws.Range("P6").Value = 0.08772724768620538
$ws.Range("Q6").Value = 0.04518009808644445
$ws.Range("R6").Value = 0.406620882778
$ws.Range("S6").Value = 0.04679920318699909
$ws.Range("T6").Value = 0.04679920318699908

# Row 7
$ws.Range("G7").Value = 2.127891333333333
$ws.Range("H7").Value = 6.383674
$ws.Range("I7").Value = 0.5334625720209161
$ws.Range("J7").Value = 0.5334625720209161
$ws.Range("N7").Value = 0.6623830000000001
$ws.Range("O7").Value = 0.9122727523137947
$ws.Range("P7").Value = 0.9122727523137947
$ws.Range("Q7").Value = 0.4698263483491111
$ws.Range("S7").Value = 0.486663368833917
$ws.Range("T7").Value = 0.486663368833917
